# Weekly update: a new observation is added for
# "Terminal La Palmera de La Serena - Papa" which shifts all the existing
# records (rows 424..442) down by one row (to 425..443) and inserts a
# brand-new row 424 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 424; everything below
# (old rows 424-442) shifts down to 425-443, preserving all their values
# and formatting automatically.
$ws.Rows.Item(424).Insert()

# Populate the newly inserted row 424 with this week's record. Most of
# the "dimension" columns repeat the same market/category as the row
# that used to be here (now row 425), only the measured fields change.
$ws.Cells.Item(424, 1).Value2  = 8
$ws.Cells.Item(424, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(424, 3).Value2  = "Coquimbo"
$ws.Cells.Item(424, 4).Value2  = 44753
$ws.Cells.Item(424, 5).Value2  = 4
$ws.Cells.Item(424, 6).Value2  = 100114001
$ws.Cells.Item(424, 7).Value2  = "Papa"
$ws.Cells.Item(424, 8).Value2  = "Asterix"
$ws.Cells.Item(424, 9).Value2  = "1a (guarda)"
$ws.Cells.Item(424, 10).Value2 = 2400
$ws.Cells.Item(424, 11).Value2 = 9500
$ws.Cells.Item(424, 12).Value2 = 10000
$ws.Cells.Item(424, 13).Value2 = 9750
$ws.Cells.Item(424, 14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item(424, 15).Value2 = "Región de Los Lagos"
$ws.Cells.Item(424, 16).Value2 = 390
$ws.Cells.Item(424, 17).Value2 = 25
$ws.Cells.Item(424, 18).Value2 = "Hortaliza"
